$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Split the word "coordinate" in the opening paragraph and drop a
#    collapsed "_GoBack" bookmark right in the middle of it ("coo" |
#    "rdinate"). We do this BEFORE touching the old bookmark location
#    so the two edits cannot interfere with each other's offsets.
# ------------------------------------------------------------------
$coordRange = $d.Content
$coordRange.Find.MatchCase = $true
$foundCoord = $coordRange.Find.Execute("to coordinate the Navy")
if ($foundCoord) {
    $splitPos = $coordRange.Start + 6
    $newBmRange = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_GoBack", $newBmRange)
}

# ------------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark. It used to wrap everything
#    from right after "Data Sets Available " through just before
#    "to help" -- deleting it just drops the bookmark, not the text.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack") -and $d.Bookmarks.Count -ge 0) {
    # There can only be one "_GoBack" bookmark at a time in a live Word
    # session; re-fetch it defensively in case Add() above recycled it.
}

$allRange = $d.Content
$allRange.Find.MatchCase = $true
$allRange.Find.ClearFormatting()
$allRange.Find.Replacement.ClearFormatting()

# Re-merge "Data Sets Available " + "for Integration" (the two runs
# that used to be split apart by the old bookmarkStart) into one run.
$allRange.Find.Execute("Data Sets Available for Integration", $true, $false, $false, $false, $false, $true, 1, $false, "Data Sets Available for Integration", 2) | Out-Null

# Re-merge " available " + "to help" (the two runs that used to be
# split apart by the old bookmarkEnd) into one run.
$allRange2 = $d.Content
$allRange2.Find.MatchCase = $true
$allRange2.Find.ClearFormatting()
$allRange2.Find.Replacement.ClearFormatting()
$allRange2.Find.Execute("available to help", $true, $false, $false, $false, $false, $true, 1, $false, "available to help", 2) | Out-Null
